# Adding some Yoddle related API Questions. (#37)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24: fill in Screen / API / Issue for the existing "PFM" row
$ws.Range("B24").Value = "Finance and Category Screen"
$ws.Range("C24").Value = "Yoddle"
$ws.Range("F24").Value = "What all endpoints we will be using for Income and expense for a category?"

# Row 25: fill in Screen / API / Issue for the existing "PFM" row
$ws.Range("B25").Value = "Category Secreen"
$ws.Range("C25").Value = "Yoddle"
$ws.Range("F25").Value = "What all endpoints we will be using for first level category on category Screen"

# Row 26: brand new row continuing the Stream/Screen/API/Issue pattern
$ws.Range("A26").Value = "PFM"
$ws.Range("B26").Value = "Transaction"
$ws.Range("C26").Value = "Yoddle"
$ws.Range("F26").Value = "What all endpoints we will be using for second level category on category Screen"

# Widen column F slightly to fit the new, longer issue text
$ws.Columns.Item(6).ColumnWidth = 76.8

# Match the author's final selection/scroll position
$excel.Goto($ws.Range("A19"), $true)
$ws.Range("B26").Select()
